$d = $word.ActiveDocument

# Locate the "postal address:" term paragraph and remove it together with
# the immediately following "43 Koerselsebaan ... Belgium" definition
# paragraph (the two form a single definition-list entry that was dropped).
$find = $d.Content
$find.Find.ClearFormatting()
$found = $find.Find.Execute("postal address:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $termPara = $find.Paragraphs(1)
    $defPara = $termPara.Next()

    $start = $termPara.Range.Start
    $end = $defPara.Range.End

    $d.Range($start, $end).Delete()
}
